$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: B2, C2 updated; E2 cleared
$ws.Range("B2").Value = "NSE:NDRAUTO"
$ws.Range("C2").Value = "NSE:DCI"
$ws.Range("E2").Value = ""

# Row 3: B3 cleared; C3 updated; E3 cleared
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = "NSE:FIVESTAR"
$ws.Range("E3").Value = ""

# Row 4: B4 cleared; C4 updated; E4 cleared
$ws.Range("B4").Value = ""
$ws.Range("C4").Value = "NSE:HOVS"
$ws.Range("E4").Value = ""

# Row 5: B5 cleared; C5 updated; E5 cleared
$ws.Range("B5").Value = ""
$ws.Range("C5").Value = "NSE:LYPSAGEMS"
$ws.Range("E5").Value = ""

# Row 6: B6 cleared; C6 updated; E6 cleared
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = "NSE:NILKAMAL"
$ws.Range("E6").Value = ""

# Row 7: B7 cleared; C7 updated; E7 cleared
$ws.Range("B7").Value = ""
$ws.Range("C7").Value = "NSE:PRIMESECU"
$ws.Range("E7").Value = ""

# Remove rows 8 through 14 entirely
$ws.Rows("8:14").Delete()
